$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# New "Execute" / "TestCase" columns (A & B) with the data-provider
# rows, written in the same order the shared-string table records
# them in (header row B-then-A, then column B down, then column A
# down) so the unique-string ordering matches exactly.
# ------------------------------------------------------------------
$ws.Range("B1").Value = "TestCase"
$ws.Range("A1").Value = "Execute"

$ws.Range("B2").Value = "TC001"
$ws.Range("B3").Value = "TC002"
$ws.Range("B4").Value = "TC003"
$ws.Range("B5").Value = "TC004"

$ws.Range("A2").Value = "Y"
$ws.Range("A3").Value = "N"
$ws.Range("A4").Value = "Y"
$ws.Range("A5").Value = "Y"

# ------------------------------------------------------------------
# Existing Username / Password / Error table, shifted right from
# A:C to C:E, plus two new rows and the numeric Error-code column.
# ------------------------------------------------------------------
$ws.Range("C1").Value = "Username"
$ws.Range("D1").Value = "Password"
$ws.Range("E1").Value = "Error"

$ws.Range("C2").Value = "somesh"
$ws.Range("D2").Value = "test"
$ws.Range("E2").Value = 1111

$ws.Range("C3").Value = "test"
$ws.Range("D3").Value = "admin"
$ws.Range("E3").Value = 2222

$ws.Range("C4").Value = "admin"
$ws.Range("D4").Value = "test"
$ws.Range("E4").Value = 3333

$ws.Range("C5").Value = "admin"
$ws.Range("D5").Value = "admin"
$ws.Range("E5").Value = 4444

# ------------------------------------------------------------------
# Best-fit-style column widths left behind on the Username/Password
# columns (C & D).
# ------------------------------------------------------------------
$ws.Columns.Item(3).ColumnWidth = 9.15
$ws.Columns.Item(4).ColumnWidth = 8.65

# ------------------------------------------------------------------
# Selection ends up on E9.
# ------------------------------------------------------------------
$ws.Range("E9").Select() | Out-Null

Write-Output "done"
